$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(5)
$ws.Range("D1:D12").EntireColumn.AutoFit()
